$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.Formula = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '63.821.57'
$ws.Range('E2').Value = '  +0.78%  '
Set-TextValue $ws 'D3' '3.319.69'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue $ws 'D5' '604.67'
$ws.Range('E5').Value = '  +1.91%  '
Set-TextValue $ws 'D6' '142.70'
$ws.Range('E6').Value = '  +0.87%  '
Set-TextValue $ws 'D8' '3.319.42'
$ws.Range('E8').Value = '  +2.65%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +1.72%  '
Set-TextValue $ws 'D11' '5.54'
$ws.Range('E11').Value = '  +3.77%  '
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('E13').Value = '  +0.47%  '
Set-TextValue $ws 'D14' '35.03'
$ws.Range('E14').Value = '  +2.09%  '
Set-TextValue $ws 'D15' '3.872.42'
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('E16').Value = '  +0.34%  '
Set-TextValue $ws 'D17' '3.325.67'
$ws.Range('E17').Value = '  +2.60%  '
Set-TextValue $ws 'D18' '63.934.20'
$ws.Range('E18').Value = '  +0.95%  '
Set-TextValue $ws 'D19' '6.86'
$ws.Range('E19').Value = '  +1.42%  '
Set-TextValue $ws 'D20' '481.36'
$ws.Range('E20').Value = '  +1.22%  '
Set-TextValue $ws 'D21' '14.09'
$ws.Range('E21').Value = '  -0.05%  '
Set-TextValue $ws 'D22' '0.739'
$ws.Range('E22').Value = '  +2.19%  '
Set-TextValue $ws 'D23' '7.97'
$ws.Range('E23').Value = '  +0.92%  '
Set-TextValue $ws 'D24' '13.94'
$ws.Range('E24').Value = '  +5.83%  '
Set-TextValue $ws 'D25' '84.94'
$ws.Range('E25').Value = '  +1.11%  '
Set-TextValue $ws 'D27' '2.78'
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D28' '1.00'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D29' '8.26'
$ws.Range('E29').Value = '  +2.81%  '
Set-TextValue $ws 'D30' '7.17'
$ws.Range('E30').Value = '  -3.89%  '
$ws.Range('E31').Value = '  +2.63%  '
Set-TextValue $ws 'D32' '28.88'
$ws.Range('E32').Value = '  +5.05%  '
Set-TextValue $ws 'D33' '0.106'
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('E35').Value = '  +1.36%  '
$ws.Range('E36').Value = '  +3.31%  '
Set-TextValue $ws 'D37' '0.0₃0748'
$ws.Range('E37').Value = '  +5.25%  '
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D40' '3.132.71'
$ws.Range('E40').Value = '  +5.27%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 'D41' '434.09'
$ws.Range('E41').Value = '  +2.99%  '
$ws.Range('E42').Value = '  +7.11%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 'D43' '8.36'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D44' '2.76'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('E45').Value = '  +0.24%  '
Set-TextValue $ws 'D46' '2.25'
$ws.Range('E46').Value = '  +4.17%  '
Set-TextValue $ws 'D47' '36.76'
$ws.Range('E47').Value = '  +9.05%  '
Set-TextValue $ws 'D48' '26.42'
$ws.Range('E48').Value = '  +2.32%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D51' '0.114'
$ws.Range('E51').Value = '  -0.62%  '
